# Updated cryptos list values (Price / Volume(1h)) for rows 2-51 of the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new text looks like a plain number (e.g. "551.62"). Assigning
# such a string directly would make Excel auto-convert the cell to a Number, which
# would not match the original inline-string "Price" column. Force these to remain
# text the same way a user would in the UI: type them with a leading apostrophe
# (quote-prefix), then clear the resulting cell formatting so the style reverts to
# the default (unstyled) look used throughout the rest of the column.
$numericLookingPriceRows = @(5,6,7,8,9,12,13,20,23,25,30,31,34,39,40,41,43,45,46,49,50,51)

$ws.Range("D2").Value = '59.776.10'
$ws.Range("E2").Value = '  +2.62%  '
$ws.Range("D3").Value = '2.419.28'
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''551.62'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").Value = '''137.33'
$ws.Range("E6").Value = '  +2.92%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''0.587'
$ws.Range("E8").Value = '  +3.09%  '
$ws.Range("D9").Value = '''0.106'
$ws.Range("E9").Value = '  -0.62%  '
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("E11").Value = '  -2.23%  '
$ws.Range("D12").Value = '''0.355'
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").Value = '''24.93'
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("D14").Value = '2.847.35'
$ws.Range("E14").Value = '  +2.24%  '
$ws.Range("D15").Value = '59.724.99'
$ws.Range("E15").Value = '  +2.75%  '
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = '2.413.93'
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("D20").Value = '''329.42'
$ws.Range("E20").Value = '  -0.63%  '
$ws.Range("E21").Value = '  -3.67%  '
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '''66.04'
$ws.Range("E23").Value = '  +3.42%  '
$ws.Range("E24").Value = '  +1.15%  '
$ws.Range("D25").Value = '''8.82'
$ws.Range("E25").Value = '  +6.34%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  +3.24%  '
$ws.Range("D28").Value = '0.0₃0774'
$ws.Range("E28").Value = '  +4.28%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").Value = '''170.04'
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("D31").Value = '''6.11'
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("D34").Value = '''0.999'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  +3.54%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("D39").Value = '''39.59'
$ws.Range("E39").Value = '  -2.05%  '
$ws.Range("D40").Value = '''0.409'
$ws.Range("E40").Value = '  -6.57%  '
$ws.Range("D41").Value = '''312.93'
$ws.Range("E41").Value = '  +8.44%  '
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("D43").Value = '''138.75'
$ws.Range("E43").Value = '  -1.79%  '
$ws.Range("E44").Value = '  +1.20%  '
$ws.Range("D45").Value = '''0.0518'
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").Value = '''19.43'
$ws.Range("E46").Value = '  +4.40%  '
$ws.Range("E47").Value = '  +1.91%  '
$ws.Range("E48").Value = '  +0.43%  '
$ws.Range("D49").Value = '''0.392'
$ws.Range("E49").Value = '  -6.17%  '
$ws.Range("D50").Value = '''17.57'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '''11.06'
$ws.Range("E51").Value = '  +0.22%  '

foreach ($r in $numericLookingPriceRows) {
    $ws.Range("D$r").ClearFormats()
}

Write-Host "Updated cryptos list"
